$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string text content (row 3: internet -> water) ---
$ws.Range("A3").Value = "water"
$ws.Range("B3").Value = "Water"
$ws.Range("C3").Value = "Last technology pipes, granted heat water even during the winter "

# --- row 4: garbage -> hydro ---
$ws.Range("A4").Value = "hydro"
$ws.Range("B4").Value = "Hydro"
$ws.Range("C4").Value = "Granted electric power supply 365 days a year "

# --- row 5: electricity description updated (only B5 text changes) ---
$ws.Range("B5").Value = "Electricity"

# --- row 6: television -> gas ---
$ws.Range("A6").Value = "gas"
$ws.Range("B6").Value = "Gas"
$ws.Range("C6").Value = "Gas burners included in the kitchen and balcony "

# --- row 7: telephone -> airConditioning ---
$ws.Range("A7").Value = "airConditioning"
$ws.Range("B7").Value = "Air Conditioning"
$ws.Range("C7").Value = "Great savings for your summer bills!! Intelligent system of air conditioning in all rooms "

# --- Remove underline formatting from the C6 cell's font ---
$ws.Range("C6").Font.Underline = -4142

# --- Column widths: split B:C into its own widths (column C keeps its original width) ---
$ws.Columns.Item(2).ColumnWidth = 22.67

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 49.25
$ws.Rows.Item(6).RowHeight = 49.25
$ws.Rows.Item(7).RowHeight = 73.1

# --- Selection moves from B4 to C5 ---
[void]$ws.Range("C5").Select()

# --- Tab ratio (window split position) ---
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.203
